$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 ("5816812 - João Paulo Alves Silva" with no label in column A)
# was removed; everything below it shifted up one row.
$ws.Rows(13).Delete()

# After the shift, a handful of long-text value cells were replaced with new
# (shorter) content while keeping their row's label in column A.
$ws.Range("B10").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C10").Value = "5816812 - João Paulo Alves Silva"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

$ws.Range("B18").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C18").Value = "5816812 - João Paulo Alves Silva"
